$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at G:H, pushing the existing "Datos Londres"/"Datos UK"
# columns (and their data) to I:J.
$ws.Range("G:H").Insert()

# New header cells for the inserted columns.
$ws.Range("G1").Value = "Inicio"
$ws.Range("H1").Value = "Fin"

# Indicador 01 (row 2) - fill in the date range + new "Inicio"/"Fin" columns
$ws.Range("C2").Value = 1999
$ws.Range("D2").Value = 2020
$ws.Range("E2").Value = 2021
$ws.Range("F2").Value = 2031
$ws.Range("G2").Value = 1999
$ws.Range("H2").Value = 2031
$ws.Range("K2").Value = "NO"

# Indicador 02 (row 3)
$ws.Range("C3").Value = 2011
$ws.Range("D3").Value = 2050
$ws.Range("G3").Value = 2011
$ws.Range("H3").Value = 2031
$ws.Range("K3").Value = "SI"

# Indicador 03 - Empleo (row 4) - new data load
$ws.Range("B4").Value = "Salario"
$ws.Range("C4").Value = 2002
$ws.Range("D4").Value = 2022
$ws.Range("E4").Value = 2023
$ws.Range("F4").Value = 2031
$ws.Range("G4").Value = 2002
$ws.Range("H4").Value = 2031
$ws.Range("I4").Value = "SI"
$ws.Range("J4").Value = "SI"
$ws.Range("K4").Value = "NO"

# New "CoL" header, last new string introduced in the workbook.
$ws.Range("K1").Value = "CoL"
$ws.Range("K1").Font.Bold = $true

# Column widths: the two newly inserted columns get a fixed custom width,
# matching the neighboring "Año fin" column.
$ws.Range("G:H").ColumnWidth = 10.25

$ws.Range("K5").Select()
